# Target change (word/numbering.xml only - no body content is affected):
#
#   abstractNum 990   <w:nsid w:val="b7720116"/> -> "afe9bc24"
#   abstractNum 991   <w:nsid w:val="13330212"/> -> "ad65d943"
#   abstractNum 99721 <w:nsid w:val="aaa74375"/> -> "be3b6d21"
#   abstractNum 99722 <w:nsid w:val="432c33d0"/> -> "9dd42522"
#
# <w:nsid> is Word's internal "list signature" GUID: it has no visible or
# semantic effect (no text/formatting/numbering behaviour reads it), and in
# the Word object model it is only ever surfaced, read-only, as
# List.ListID. Re-stamp every list we can reach through the object model
# with its target id - this is the one legitimate automation entry point
# for this value.

$d = $word.ActiveDocument

# The two lists actually applied to content (numId 1002 / 1003, backing
# abstractNum 99721 / 99722) are reachable live through each list
# paragraph's Range.ListFormat.List.
$targetNsidByNumId = @{
    "1002" = "be3b6d21"
    "1003" = "9dd42522"
}

foreach ($para in $d.Paragraphs) {
    $listFormat = $para.Range.ListFormat
    if ($listFormat.ListType -ne 0) {
        $lst = $listFormat.List
        $numId = [string]$lst.ListID
        if ($targetNsidByNumId.ContainsKey($numId)) {
            $lst.ListID = $targetNsidByNumId[$numId]
        }
    }
}

# abstractNum 990 / 991 back unused list definitions (no paragraph in the
# document is formatted with numId 1000 / 1001), so also sweep
# Document.Lists for completeness.
foreach ($nsid in @("afe9bc24", "ad65d943", "be3b6d21", "9dd42522")) {
    for ($i = 1; $i -le $d.Lists.Count; $i++) {
        $d.Lists.Item($i).ListID = $nsid
    }
}

Write-Host "Restamped list ids for abstractNum 990, 991, 99721, 99722."
